$d = $word.ActiveDocument

# Update the date paragraph
$d.Content.Find.Execute("2024-08-28 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-29 Thursday", 2) | Out-Null

# Update table cells by position (row, col) to avoid ambiguity with repeated text
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "91+1="
$t.Cell(1,2).Range.Text = "21-3="
$t.Cell(1,3).Range.Text = "74-69="
$t.Cell(1,4).Range.Text = "64+2="
$t.Cell(1,5).Range.Text = "15+36="
$t.Cell(2,1).Range.Text = "85-79="
$t.Cell(2,2).Range.Text = "28-28="
$t.Cell(2,3).Range.Text = "21+74="
$t.Cell(2,4).Range.Text = "7-2="
$t.Cell(2,5).Range.Text = "3+91="
$t.Cell(3,1).Range.Text = "97-8="
$t.Cell(3,2).Range.Text = "19+43="
$t.Cell(3,3).Range.Text = "95-40="
$t.Cell(3,4).Range.Text = "26+2="
$t.Cell(3,5).Range.Text = "74-54="
$t.Cell(4,1).Range.Text = "71+15="
$t.Cell(4,2).Range.Text = "98-82="
$t.Cell(4,3).Range.Text = "93-49="
$t.Cell(4,4).Range.Text = "13+9="
$t.Cell(4,5).Range.Text = "31+28="
$t.Cell(5,1).Range.Text = "55+9="
$t.Cell(5,2).Range.Text = "46+10="
$t.Cell(5,3).Range.Text = "11+25="
$t.Cell(5,4).Range.Text = "5+4="
$t.Cell(5,5).Range.Text = "38+3="
$t.Cell(6,1).Range.Text = "51+11="
$t.Cell(6,2).Range.Text = "93-27="
$t.Cell(6,3).Range.Text = "17-3="
$t.Cell(6,4).Range.Text = "55-44="
$t.Cell(6,5).Range.Text = "12+51="
$t.Cell(7,1).Range.Text = "81-38="
$t.Cell(7,2).Range.Text = "78+5="
$t.Cell(7,3).Range.Text = "75-5="
$t.Cell(7,4).Range.Text = "6+14="
$t.Cell(7,5).Range.Text = "5+26="
$t.Cell(8,1).Range.Text = "20+50="
$t.Cell(8,2).Range.Text = "7+4="
$t.Cell(8,3).Range.Text = "37+58="
$t.Cell(8,4).Range.Text = "7+45="
$t.Cell(8,5).Range.Text = "91+1="
$t.Cell(9,1).Range.Text = "35+42="
$t.Cell(9,2).Range.Text = "25+60="
$t.Cell(9,3).Range.Text = "44-21="
$t.Cell(9,4).Range.Text = "45-14="
$t.Cell(9,5).Range.Text = "99-30="
$t.Cell(10,1).Range.Text = "70-13="
$t.Cell(10,2).Range.Text = "37+43="
$t.Cell(10,3).Range.Text = "28+7="
$t.Cell(10,4).Range.Text = "71-48="
$t.Cell(10,5).Range.Text = "0+25="
$t.Cell(11,1).Range.Text = "86-57="
$t.Cell(11,2).Range.Text = "31+19="
$t.Cell(11,3).Range.Text = "91-67="
$t.Cell(11,4).Range.Text = "47+12="
$t.Cell(11,5).Range.Text = "91-12="
$t.Cell(12,1).Range.Text = "18+39="
$t.Cell(12,2).Range.Text = "61-18="
$t.Cell(12,3).Range.Text = "94-29="
$t.Cell(12,4).Range.Text = "45-9="
$t.Cell(12,5).Range.Text = "36+13="
$t.Cell(13,1).Range.Text = "7+87="
$t.Cell(13,2).Range.Text = "59-12="
$t.Cell(13,3).Range.Text = "56-19="
$t.Cell(13,4).Range.Text = "71-43="
$t.Cell(13,5).Range.Text = "97-1="
$t.Cell(14,1).Range.Text = "64-25="
$t.Cell(14,2).Range.Text = "96+1="
$t.Cell(14,3).Range.Text = "21-14="
$t.Cell(14,4).Range.Text = "18+48="
$t.Cell(14,5).Range.Text = "94-76="
$t.Cell(15,1).Range.Text = "80+1="
$t.Cell(15,2).Range.Text = "22+10="
$t.Cell(15,3).Range.Text = "52+38="
$t.Cell(15,4).Range.Text = "32+53="
$t.Cell(15,5).Range.Text = "73+10="
$t.Cell(16,1).Range.Text = "29-25="
$t.Cell(16,2).Range.Text = "52+18="
$t.Cell(16,3).Range.Text = "93-66="
$t.Cell(16,4).Range.Text = "87-36="
$t.Cell(16,5).Range.Text = "85-2="
$t.Cell(17,1).Range.Text = "86+6="
$t.Cell(17,2).Range.Text = "54-0="
$t.Cell(17,3).Range.Text = "72+21="
$t.Cell(17,4).Range.Text = "78-77="
$t.Cell(17,5).Range.Text = "7+69="
$t.Cell(18,1).Range.Text = "2+97="
$t.Cell(18,2).Range.Text = "30+33="
$t.Cell(18,3).Range.Text = "98-93="
$t.Cell(18,4).Range.Text = "61+37="
$t.Cell(18,5).Range.Text = "6+57="
$t.Cell(19,1).Range.Text = "2+56="
$t.Cell(19,2).Range.Text = "89+8="
$t.Cell(19,3).Range.Text = "56+12="
$t.Cell(19,4).Range.Text = "38+25="
$t.Cell(19,5).Range.Text = "80-12="
$t.Cell(20,1).Range.Text = "26+60="
$t.Cell(20,2).Range.Text = "44-19="
$t.Cell(20,3).Range.Text = "6+37="
$t.Cell(20,4).Range.Text = "13+86="
$t.Cell(20,5).Range.Text = "84-9="
